# Weekly update: insert a new price observation as the new first data row
# (row 8) for "Vega Monumental Concepción - Haba", pushing all subsequent
# rows down by one. This mirrors the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8; existing rows 8..63 shift to 9..64.
$ws.Rows("8:8").Insert()

# Populate the newly inserted row 8 with the new weekly record.
$ws.Range("A8").Value = 11
$ws.Range("B8").Value = "Vega Monumental Concepción"
$ws.Range("C8").Value = "Bíobío"
$ws.Range("D8").Value = 45163
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 100112026
$ws.Range("G8").Value = "Haba"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 50
$ws.Range("K8").Value = 13000
$ws.Range("L8").Value = 13000
$ws.Range("M8").Value = 13000
$ws.Range("N8").Value = "`$/saco 25 kilos"
$ws.Range("O8").Value = "Región de Coquimbo"
$ws.Range("P8").Value = 520
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of column D.
$ws.Range("D8").NumberFormat = $ws.Range("D9").NumberFormat
